# Employee Import.xlsx — close the AN:AT header gap and append a new
# "IsNonProduction" header column.
#
# Before: headers sit in A1:AM1 (cols 1-39), then columns AN1:AT1 (40-46)
#         are blank, then headers resume in AU1:BZ1 (47-78).
# After:  headers sit in A1:AM1 (cols 1-39) followed *contiguously* by
#         AN1:BT1 (40-72) -- i.e. every header that used to start at AU1
#         shifts left by 7 columns into the old gap -- and a brand new
#         header "IsNonProduction" is appended as the new last column
#         (BT1, col 72). Columns BU1:BZ1, now past the used range, are
#         cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 32 header values that used to live in AU1:BZ1 (cols 47-78), in
# order, plus the one brand new header appended at the end.
$vals = @(
    "MiddleName",
    "OfficialPhone",
    "PersonalLocation",
    "PersonalEmail",
    "LeaveGroup",
    "CompanyMaster",
    "LocationMaster",
    "HolidayCalendar",
    "Status",
    "AadharNo",
    "PanNo",
    "PassportNo",
    "DrivingLicense",
    "BankName",
    "BankAccountNo",
    "BankIfscCode",
    "BankBranch",
    "HomeAddress",
    "FatherName",
    "EmergencyContactPerson1",
    "EmergencyContactPerson2",
    "EmergencyContactNo1",
    "EmergencyContactNo2",
    "MotherName",
    "FatherAadharNo",
    "MotherAadharNo",
    "OrganizationType",
    "WorkingStatus",
    "ConfirmationDate",
    "PostalCode",
    "ApprovalLevel",
    "OfficialEmail",
    "IsNonProduction"
)

# Write them starting at column 40 (AN) through column 72 (BT).
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(1, 40 + $i).Value = $vals[$i]
}

# The old tail (BU1:BZ1, cols 73-78) is now outside the used range --
# clear it out so it doesn't linger as stale data.
$null = $ws.Range("BU1:BZ1").ClearContents()

# Match the saved selection/viewport: AN1:BT1 selected with AN1 active.
$null = $ws.Range("AN1:BT1").Select()
